$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# "Experimental" row: the Value cell (B7) needs to become the literal text "true"
# (not the native Boolean TRUE that a plain .Value assignment would produce).
# Route it through a formula -> copy -> paste-values round trip so it lands
# in the sheet as a text cell.
$ws.Range("B7").Formula = "=""true"""
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)

# "Date" row: update the timestamp text value.
$ws.Range("B8").Value = "2023-02-16T14:43:10-06:00"
